# Read External Data from Excel
# - Expand the short names in A2/A3 into full names.
# - Turn on word-wrap for those two cells so the longer names display
#   nicely, which mints a new cell style (wrapText alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sebastian Baron Caicedo"
$ws.Range("A3").Value = "Luis Carlos Covilla Yarce"

$ws.Range("A2:A3").WrapText = $true
